$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet "Overview"
# =====================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() on this engine clears the whole sheet collection;
# rebuild everything afterwards in final left-to-right, top-to-bottom order.
$ws1.Range("A1").Hyperlinks.Delete()

# Push the ".localization-config" row from row 4 down to row 6, and fill in
# the two new "Ready for handoff" rows at 4 and 5.
$ws1.Range("A6").Value = ".localization-config"
$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"

$ws1.Range("A4").Value = "099a1c65-a65c-4a8e-b088-67a3a96283e3.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = "3338663b-41cd-4af3-8cd2-89dee62ff182.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/5b69114e-d216-4c25-9281-a397e47b2e6e.md", "", "", "5b69114e-d216-4c25-9281-a397e47b2e6e.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/d32161cc-5585-4806-b1a0-df8c2e7ba787.md", "", "", "d32161cc-5585-4806-b1a0-df8c2e7ba787.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/099a1c65-a65c-4a8e-b088-67a3a96283e3.md", "", "", "099a1c65-a65c-4a8e-b088-67a3a96283e3.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/3338663b-41cd-4af3-8cd2-89dee62ff182.md", "", "", "3338663b-41cd-4af3-8cd2-89dee62ff182.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/.localization-config", "", "", ".localization-config")

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

# Push the ".localization-config" row from row 4 down to row 6.
$ws2.Range("A6").Value = ".localization-config"
$ws2.Range("B6").Value = "Not to be localized"
$ws2.Range("D6").Value = "0001-01-01 00:00:00"
$ws2.Range("G6").Value = "0001-01-01 00:00:00"
$ws2.Range("H6").Value = "Ignored"

# New row 4: 099a1c65...
$ws2.Range("A4").Value = "099a1c65-a65c-4a8e-b088-67a3a96283e3.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-09 12:39:10"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Include"

# New row 5: 3338663b...
$ws2.Range("A5").Value = "3338663b-41cd-4af3-8cd2-89dee62ff182.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.zh-cn.xlf"
$ws2.Range("D5").Value = "2016-03-09 12:39:10"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/5b69114e-d216-4c25-9281-a397e47b2e6e.md", "", "", "5b69114e-d216-4c25-9281-a397e47b2e6e.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4e25e85976451c4da9497bb10a7ebc4df8f03f8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5b69114e-d216-4c25-9281-a397e47b2e6e.62df8269a6d767521d87c7fc4a12ba91fbda4f62.zh-cn.xlf", "", "", "5b69114e-d216-4c25-9281-a397e47b2e6e.62df8269a6d767521d87c7fc4a12ba91fbda4f62.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/d32161cc-5585-4806-b1a0-df8c2e7ba787.md", "", "", "d32161cc-5585-4806-b1a0-df8c2e7ba787.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4e25e85976451c4da9497bb10a7ebc4df8f03f8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d32161cc-5585-4806-b1a0-df8c2e7ba787.0661a84c2ea57c5b79b58498a04a5f3d3602d9b8.zh-cn.xlf", "", "", "d32161cc-5585-4806-b1a0-df8c2e7ba787.0661a84c2ea57c5b79b58498a04a5f3d3602d9b8.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/099a1c65-a65c-4a8e-b088-67a3a96283e3.md", "", "", "099a1c65-a65c-4a8e-b088-67a3a96283e3.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4e25e85976451c4da9497bb10a7ebc4df8f03f8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.zh-cn.xlf", "", "", "099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/3338663b-41cd-4af3-8cd2-89dee62ff182.md", "", "", "3338663b-41cd-4af3-8cd2-89dee62ff182.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4e25e85976451c4da9497bb10a7ebc4df8f03f8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.zh-cn.xlf", "", "", "3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/.localization-config", "", "", ".localization-config")

# =====================================================================
# Sheet "de-de"
# =====================================================================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

# Push the ".localization-config" row from row 4 down to row 6.
$ws3.Range("A6").Value = ".localization-config"
$ws3.Range("B6").Value = "Not to be localized"
$ws3.Range("D6").Value = "0001-01-01 00:00:00"
$ws3.Range("G6").Value = "0001-01-01 00:00:00"
$ws3.Range("H6").Value = "Ignored"

# New row 4: 099a1c65...
$ws3.Range("A4").Value = "099a1c65-a65c-4a8e-b088-67a3a96283e3.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-09 12:39:20"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Include"

# New row 5: 3338663b...
$ws3.Range("A5").Value = "3338663b-41cd-4af3-8cd2-89dee62ff182.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-09 12:39:20"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/5b69114e-d216-4c25-9281-a397e47b2e6e.md", "", "", "5b69114e-d216-4c25-9281-a397e47b2e6e.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf783ee87320be698da155b13e5357eb7e1483ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5b69114e-d216-4c25-9281-a397e47b2e6e.62df8269a6d767521d87c7fc4a12ba91fbda4f62.de-de.xlf", "", "", "5b69114e-d216-4c25-9281-a397e47b2e6e.62df8269a6d767521d87c7fc4a12ba91fbda4f62.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/d32161cc-5585-4806-b1a0-df8c2e7ba787.md", "", "", "d32161cc-5585-4806-b1a0-df8c2e7ba787.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf783ee87320be698da155b13e5357eb7e1483ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d32161cc-5585-4806-b1a0-df8c2e7ba787.0661a84c2ea57c5b79b58498a04a5f3d3602d9b8.de-de.xlf", "", "", "d32161cc-5585-4806-b1a0-df8c2e7ba787.0661a84c2ea57c5b79b58498a04a5f3d3602d9b8.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/099a1c65-a65c-4a8e-b088-67a3a96283e3.md", "", "", "099a1c65-a65c-4a8e-b088-67a3a96283e3.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf783ee87320be698da155b13e5357eb7e1483ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.de-de.xlf", "", "", "099a1c65-a65c-4a8e-b088-67a3a96283e3.691988ad891a3b521a6d2ad79735a4a094dc3c48.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/e2e/3338663b-41cd-4af3-8cd2-89dee62ff182.md", "", "", "3338663b-41cd-4af3-8cd2-89dee62ff182.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bf783ee87320be698da155b13e5357eb7e1483ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.de-de.xlf", "", "", "3338663b-41cd-4af3-8cd2-89dee62ff182.aea0f8ce1f0ce3aa5b52729717075716ae5de333.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/cb32449fec2295264de40191640eb58418e047ea/.localization-config", "", "", ".localization-config")
